$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.405.49"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.54%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.287.13"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.92%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "156.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +15,584.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "307.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.96%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "95.62"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.72%  "
$ws.Range("E8").Value = "  +0.24%  "
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("E10").Value = "  +2.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "35.90"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +11.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0803"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.09%  "
$ws.Range("E13").Value = "  -1.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.51%  "
$ws.Range("E15").Value = "  +0.98%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.47"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.292.41"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.797"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.307.29"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.65"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0918"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.74%  "
$ws.Range("E22").Value = "  +1.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.03"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "243.04"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.54%  "
$ws.Range("E25").Value = "  +0.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.94"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.21%  "
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "36.07"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.78%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.59"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.89%  "
$ws.Range("E31").Value = "  -8.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "161.94"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.34"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.80%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0754"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.54%  "
$ws.Range("E36").Value = "  +2.73%  "
$ws.Range("E37").Value = "  +4.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.26"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.38"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.13%  "
$ws.Range("E40").Value = "  +2.35%  "
$ws.Range("E41").Value = "  -0.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.18"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.013.28"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.51"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.31%  "
$ws.Range("E45").Value = "  +10.91%  "
$ws.Range("E46").Value = "  +2.55%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.15"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.98"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.23%  "
$ws.Range("E49").Value = "  +2.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.39"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.97"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.10%  "
